# Agregue la interfaz y el metodo para agregar usuario
# Inserts a header row on top of the employee table and sets up the
# "Numero de Nomina" (payroll number) column to auto-increment via a
# formula, plus bumps two salaries that were corrected by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new header row at the top; existing rows 1-20 (data) shift
#    down to rows 2-21.
$ws.Rows.Item(1).Insert()

# 2) Populate the new header row.
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Apellido"
$ws.Range("C1").Value = "Puesto"
$ws.Range("D1").Value = "Sueldo Mensual"
$ws.Range("E1").Value = "Fecha de Ingreso"
$ws.Range("F1").Value = "Numero de Nomina"

# 3) Give the new header columns a sensible best-fit-like width (matches
#    the width Excel would have computed for the header text).
$ws.Columns.Item(4).ColumnWidth = 13.15
$ws.Columns.Item(5).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 16.3

# 4) Apply a 2-decimal number format to the whole "Sueldo Mensual" column
#    (now rows 2-21).
$ws.Range("D2:D21").NumberFormat = "0.00"

# 5) Correct the first two salaries by hand.
$ws.Range("D2").Value = 15000
$ws.Range("D3").Value = 16000

# 6) Renumber the "Numero de Nomina" column: row 2 becomes the new
#    starting literal, and every row below derives from the one above via
#    a formula (F3 is its own formula, F4:F21 share one formula).
$ws.Range("F2").Value = 2343001
$ws.Range("F3").Formula = "=F2+1"
$ws.Range("F4:F21").Formula = "=F3+1"

# 7) Leave the selection where the user last clicked.
$ws.Range("D3").Select()
